$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(9).Insert()
$ws.Range("I1").Value = "Distribution channel code"
$ws.Range("I2").Value = "TR"
$ws.Range("I3").Value = "GO"
$ws.Columns.Item(9).ColumnWidth = 21.7
$null = $ws.Range("I1:I1048576").Select()
